$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Starbucks (row 22) previously showed "Open" for every day (Mon-Sun);
# update to show actual opening hours instead, matching the format used
# by other rows (e.g. Subway in row 23).
$ws.Range("F22:L22").Value = "0800-2200"
